$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.884.34"
$ws.Range("D3").Value = "2.442.11"
$ws.Range("E3").Value = "  -9.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "538.37"
$ws.Range("E5").Value = "  -2.84%  "
$ws.Range("D6").Value = "146.71"
$ws.Range("E6").Value = "  -7.11%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -2.94%  "
$ws.Range("D9").Value = "2.455.69"
$ws.Range("E9").Value = "  -8.62%  "
$ws.Range("D10").Value = "'0.0990"
$ws.Range("E10").Value = "  -6.69%  "
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("E13").Value = "  -5.06%  "
$ws.Range("D14").Value = "2.880.36"
$ws.Range("E14").Value = "  -8.76%  "
$ws.Range("D15").Value = "23.91"
$ws.Range("E15").Value = "  -10.07%  "
$ws.Range("D16").Value = "58.788.40"
$ws.Range("E16").Value = "  -6.38%  "
$ws.Range("E17").Value = "  -6.24%  "
$ws.Range("D18").Value = "2.503.05"
$ws.Range("E18").Value = "  -6.66%  "
$ws.Range("D19").Value = "11.12"
$ws.Range("E19").Value = "  -6.41%  "
$ws.Range("E20").Value = "  -5.65%  "
$ws.Range("D21").Value = "323.15"
$ws.Range("E21").Value = "  -6.41%  "
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  -8.96%  "
$ws.Range("D24").Value = "60.68"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("E25").Value = "  -11.16%  "
$ws.Range("E26").Value = "  -5.16%  "
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  -6.48%  "
$ws.Range("E29").Value = "  -6.18%  "
$ws.Range("D30").Value = "0.0₃0768"
$ws.Range("E30").Value = "  -10.47%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "6.64"
$ws.Range("E31").Value = "  -8.60%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  -13.07%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "156.33"
$ws.Range("E34").Value = "  -4.71%  "
$ws.Range("D35").Value = "1.38"
$ws.Range("E35").Value = "  -6.80%  "
$ws.Range("D36").Value = "'18.40"
$ws.Range("E36").Value = "  -5.54%  "
$ws.Range("E37").Value = "  -9.64%  "
$ws.Range("E38").Value = "  -5.22%  "
$ws.Range("D39").Value = "5.85"
$ws.Range("E39").Value = "  -6.46%  "
$ws.Range("D40").Value = "313.85"
$ws.Range("E40").Value = "  -10.31%  "
$ws.Range("E41").Value = "  -5.74%  "
$ws.Range("D42").Value = "0.833"
$ws.Range("E42").Value = "  -11.94%  "
$ws.Range("D43").Value = "'3.70"
$ws.Range("E43").Value = "  -7.23%  "
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("E47").Value = "  -6.05%  "
$ws.Range("D48").Value = "0.0525"
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("E49").Value = "  -5.19%  "
$ws.Range("D50").Value = "121.51"
$ws.Range("E50").Value = "  -5.49%  "
$ws.Range("D51").Value = "18.84"
$ws.Range("E51").Value = "  -9.72%  "
